$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Period row (row 6): adjust timezone-aware formatting (joda-style) instead of java Formatter
$ws.Range("B6").Value = '${from.toString("YYYY.MM.dd HH:mm:ss")+" - "+to.toString("YYYY.MM.dd HH:mm:ss")}'

# Trip data row (row 9): pass client timezone to start/end time, and switch hyperlink to https
$ws.Range("A9").Value = '${new("org.joda.time.DateTime", trip.startTime, timezone).toString("YYYY.MM.dd HH:mm:ss")}'
$ws.Range("C9").Value = '${new("org.joda.time.DateTime", trip.endTime, timezone).toString("YYYY.MM.dd HH:mm:ss")}'
$ws.Range("B9").Value = '${util.hyperlink("".format("https://www.openstreetmap.org/?mlat=%1$f&mlon=%2$f#map=16/%1$f/%2$f", trip.startLat, trip.startLon), trip.getStartAddress() == null ? "".format("%1$f°, %2$f°", trip.startLat, trip.startLon) : trip.startAddress)}'
$ws.Range("D9").Value = '${util.hyperlink("".format("https://www.openstreetmap.org/?mlat=%1$f&mlon=%2$f#map=16/%1$f/%2$f", trip.endLat, trip.endLon), trip.getEndAddress() == null ? "".format("%1$f°, %2$f°", trip.endLat, trip.endLon) : trip.endAddress)}'

# Move active selection to D9
$ws.Range("D9").Select()
